$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ACSED10"
$ws.Range("B2").Value = 2001
$ws.Range("C2").Value = "Asha Kaale"
$ws.Range("P2").Value = "PH ladies - standard"
$ws.Range("AA2").Value = "Per Day"
$ws.Range("AC2").Value = "Custom"

$ws.Range("A3").Value = "ACSED47"
$ws.Range("B3").Value = 2002
$ws.Range("C3").Value = "Sankket sonawane"
$ws.Range("P3").Value = "PH ladies - standard"
$ws.Range("AA3").Value = "Per Day"
$ws.Range("AC3").Value = "Custom"

$ws.Range("A4").Value = "ACSED46"
$ws.Range("B4").Value = 2003
$ws.Range("C4").Value = "Amaan Shaikh"
$ws.Range("P4").Value = "PH ladies - standard"
$ws.Range("AA4").Value = "Per Day"
$ws.Range("AC4").Value = "Custom"

$ws.Range("A5").Value = "ACSED11"
$ws.Range("B5").Value = 2005
$ws.Range("C5").Value = "Ajit Shinde"
$ws.Range("P5").Value = "PH Bengal Boys Supervisor"
$ws.Range("AA5").Value = "Per Day"
$ws.Range("AC5").Value = "Custom"

$ws.Range("B6").Select()
